$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "1.000") must be
# forced to Text format first, otherwise Excel would auto-convert them to
# numeric values and strip formatting such as trailing zeros.
$textForcedCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D15",
    "D17",
    "D18",
    "D19",
    "D22",
    "D25",
    "D26",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price values (Price column, D) and volume/percentage
# values (Volume(1h) column, E) exactly as scraped from coinranking.com.

$ws.Range('D2').Value = '28.518.38'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.912.99'
$ws.Range('E3').Value = '  +2.53%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '315.14'
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').Value = '0.5157'
$ws.Range('E7').Value = '  +3.19%  '
$ws.Range('D8').Value = '0.3976'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.09909'
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('D10').Value = '1.150'
$ws.Range('E10').Value = '  +2.60%  '
$ws.Range('D11').Value = '42.28'
$ws.Range('E11').Value = '  +2.33%  '
$ws.Range('D12').Value = '6.528'
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('D14').Value = '1.921.37'
$ws.Range('E14').Value = '  +3.94%  '
$ws.Range('D15').Value = '7.466'
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').Value = '0.00001140'
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').Value = '94.54'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('D19').Value = '0.06657'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('E20').Value = '  +4.87%  '
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('D22').Value = '6.305'
$ws.Range('E22').Value = '  +4.03%  '
$ws.Range('D23').Value = '28.572.58'
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('E24').Value = '  +1.19%  '
$ws.Range('D25').Value = '2.321'
$ws.Range('E25').Value = '  +2.91%  '
$ws.Range('D26').Value = '2.683'
$ws.Range('E26').Value = '  +7.98%  '
$ws.Range('D27').Value = '2.137.82'
$ws.Range('E27').Value = '  +3.56%  '
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').Value = '157.50'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').Value = '129.30'
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('E31').Value = '  +5.64%  '
$ws.Range('D32').Value = '0.1075'
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('D33').Value = '5.740'
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('D34').Value = '3.630'
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('D35').Value = '9.863'
$ws.Range('E35').Value = '  +7.84%  '
$ws.Range('D36').Value = '0.06758'
$ws.Range('E36').Value = '  -0.72%  '
$ws.Range('D37').Value = '0.02439'
$ws.Range('E37').Value = '  +2.10%  '
$ws.Range('D38').Value = '1.270'
$ws.Range('E38').Value = '  +4.87%  '
$ws.Range('D39').Value = '0.2219'
$ws.Range('E39').Value = '  +2.43%  '
$ws.Range('E40').Value = '  +2.36%  '
$ws.Range('D41').Value = '0.6478'
$ws.Range('E41').Value = '  +3.03%  '
$ws.Range('D42').Value = '5.088'
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('D43').Value = '1.186'
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('D45').Value = '13.50'
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('D47').Value = '3.764'
$ws.Range('E47').Value = '  +1.90%  '
$ws.Range('E48').Value = '  +0.81%  '
$ws.Range('D49').Value = '2.056'
$ws.Range('E49').Value = '  +3.96%  '
$ws.Range('D50').Value = '124.66'
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('D51').Value = '1.205'
$ws.Range('E51').Value = '  +1.23%  '

# Remove the temporary Text number format so the cells retain their
# original (unstyled) appearance, matching the source workbook.
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).ClearFormats()
}
